$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.437.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.590.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  +0.89%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E7").Value = "  +0.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.37"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.74%  "
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0887"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.816.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.594.80"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.530"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.444.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.65%  "
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.106"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.400.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.12%  "
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("E36").Value = "  -8.90%  "
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.542"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.811"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.983"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.725.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.40%  "
$ws.Range("E48").Value = "  +1.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("E51").Value = "  -0.99%  "

Write-Host "Updated cryptos list"